$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / Worksheets.Item(1)) and "全部类型" (sheet4 / Worksheets.Item(4))
# both list the same events (sheet4 also includes a "演出" item),
# column F = "想去人数" (want-to-go count), which increased for several events.

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 45
$ws1.Range("F4").Value = 71
$ws1.Range("F5").Value = 514
$ws1.Range("F6").Value = 1541
$ws1.Range("F8").Value = 1140
$ws1.Range("F10").Value = 216
$ws1.Range("F11").Value = 155
$ws1.Range("F12").Value = 1
$ws1.Range("F13").Value = 2
$ws1.Range("F14").Value = 1
$ws1.Range("F15").Value = 224
$ws1.Range("F16").Value = 124
$ws1.Range("F17").Value = 187
$ws1.Range("F18").Value = 174

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 45
$ws4.Range("F4").Value = 71
$ws4.Range("F5").Value = 514
$ws4.Range("F6").Value = 1541
$ws4.Range("F9").Value = 1140
$ws4.Range("F11").Value = 216
$ws4.Range("F12").Value = 155
$ws4.Range("F13").Value = 1
$ws4.Range("F14").Value = 2
$ws4.Range("F15").Value = 1
$ws4.Range("F16").Value = 225
$ws4.Range("F17").Value = 124
$ws4.Range("F18").Value = 187
$ws4.Range("F19").Value = 174
